$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.184447169303894
$ws.Range("B1").Value = 2.335700035095215
$ws.Range("C1").Value = 3.722033739089966
$ws.Range("D1").Value = 3.101845264434814
$ws.Range("E1").Value = 1.142736673355103
